$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 31 (Zapallo / Camote, "1a (cosecha)",
# origin "Región Metropolitana", dated 2023-02-21 / serial 44978). All the existing rows
# 31-55 shift down one row to 32-56.
$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44978
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112045
$ws.Range("G31").Value = "Zapallo"
$ws.Range("H31").Value = "Camote"
$ws.Range("I31").Value = "1a (cosecha)"
$ws.Range("J31").Value = 900
$ws.Range("K31").Value = 700
$ws.Range("L31").Value = 750
$ws.Range("M31").Value = 717
$ws.Range("N31").Value = "$/kilo (volumen en unidades)"
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 717
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
